$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force a run boundary at a given document character offset by
# briefly toggling Bold on the tail range and reverting it. Word's OOXML
# writer keeps runs separate once their formatting has been touched, even
# after the value is set back to the original, so this reliably produces a
# distinct <w:r> at that offset.
# ---------------------------------------------------------------------------
function Split-RunAt($offset, $endOffset) {
    $br = $d.Range($offset, $endOffset)
    $br.Bold = 1
    $br.Bold = 0
}

# ---------------------------------------------------------------------------
# 1) "Example: ..." paragraph - append a new trailing run.
# ---------------------------------------------------------------------------
$pExample = $d.Paragraphs(27)
$rExample = $pExample.Range
$rExample.MoveEnd(1, -1)
$exampleAddStart = $rExample.End
$rExample.InsertAfter(" If the trump suit had been keys, then it would" + [char]0x2019 + "ve been a matter of ")
Split-RunAt $exampleAddStart $rExample.End

# ---------------------------------------------------------------------------
# 2) "3 (Fox): ..." paragraph - append a closing parenthesis run.
# ---------------------------------------------------------------------------
$pFox = $d.Paragraphs(33)
$rFox = $pFox.Range
$rFox.MoveEnd(1, -1)
$foxAddStart = $rFox.End
$rFox.InsertAfter(")")
Split-RunAt $foxAddStart $rFox.End

# ---------------------------------------------------------------------------
# 3) "5 (Woodcutter): ..." paragraph - rewrite text and re-split into runs.
# ---------------------------------------------------------------------------
$pWood = $d.Paragraphs(34)
$rWood = $pWood.Range
$rWood.MoveEnd(1, -1)
$woodStart = $rWood.Start
$woodText = "5 (Woodcutter): When you play this, draw 1 card. Then discard any 1 card to the bottom of the deck face down (this can be the card that you just drew, if you so wish)."
$rWood.Text = $woodText

$woodParts = @(
    "5 (Woodcutter): When you play this, draw 1 card. Then discard any 1 card to t",
    "he bottom of the deck face down",
    " (t",
    "his can be the card that you just dre",
    "w, if you so wish)",
    "."
)
$offset = $woodStart
for ($i = 0; $i -lt $woodParts.Length - 1; $i++) {
    $offset = $offset + $woodParts[$i].Length
    Split-RunAt $offset $rWood.End
}

# ---------------------------------------------------------------------------
# 4) "7 (Treasure): ..." paragraph - rewrite text, re-split into runs, and
#    plant the "_GoBack" bookmark mid-word (it migrates here from the end of
#    the Monarch paragraph; Word only ever keeps a single "_GoBack" bookmark
#    so adding this one removes the old one automatically).
# ---------------------------------------------------------------------------
$pTreasure = $d.Paragraphs(35)
$rTreasure = $pTreasure.Range
$rTreasure.MoveEnd(1, -1)
$treasureStart = $rTreasure.Start
$treasureText = "7 (Treasure): After each trick, the winner receives 1 point for each 7 in the trick (meaning, at most, the winner will get 2 extra points on that one trick)."
$rTreasure.Text = $treasureText

$treasureParts = @(
    "7 (Treasure): After each trick, the winner receives ",
    "1 point for each 7 in the trick",
    " (meaning, at most, the winner will get 2 extra point",
    "s on that one trick)",
    "."
)
$offset = $treasureStart
for ($i = 0; $i -lt $treasureParts.Length - 1; $i++) {
    $offset = $offset + $treasureParts[$i].Length
    Split-RunAt $offset $rTreasure.End
}

$bookmarkOffset = $treasureStart + $treasureParts[0].Length + $treasureParts[1].Length + $treasureParts[2].Length
$bmRange = $d.Range($bookmarkOffset, $bookmarkOffset)
$d.Bookmarks.Add("_GoBack", $bmRange)
